$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B25").Value = 6472
$ws.Range("C25").Value = 1006
$ws.Range("D25").Value = 6010215
$ws.Range("E25").Value = 928.6487948084055
$ws.Range("F25").Value = 9.862502121880844
$ws.Range("G25").Value = 7.249466950959493
$ws.Range("H25").Value = 25.87243651418287
